# Deploying to gh-pages from @ Alvearie/alvearie-fhir-ig@8e4a450c507ef6f746e072652acbb72e9504f19a
# -------------------------------------------------------------------------
# Updates the "Metadata" sheet (Version/Date/Publisher/Jurisdiction, and
# drops the duplicated "Contact" row which pushes a stray last row off the
# bottom) and the "Elements" sheet (Extension's Short/Definition text).

$wb = $excel.ActiveWorkbook

# ---- Metadata sheet -------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")

# Version 5.0.0 -> 6.0.0
$meta.Range("B3").Value = "6.0.0"

# Date bumped to the new publish timestamp
$meta.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher value was blank; now populated
$meta.Range("B9").Value = "Alvearie Team"

# The old row 10 ("Contact" / "No display for ContactDetail") is replaced
# by a new "Jurisdiction" / "United States of America" row.
$meta.Range("A10").Value = "Jurisdiction"
$meta.Range("B10").Value = "United States of America"

# Row 11 duplicated the old "Contact" row entirely - remove it, which
# shifts every following row up by one (dropping the old row 21 and
# shrinking the used range to A1:B20, matching the new dimension).
$meta.Rows(11).Delete()

# ---- Elements sheet --------------------------------------------------
$elements = $wb.Worksheets.Item("Elements")

# The root Extension row's Short/Definition text changes from the generic
# "Extension" / "An Extension" to the profile-specific description.
$elements.Range("K2").Value = "Default Value"
$elements.Range("L2").Value = "Default value for the parameter"

Write-Output "edit applied"
